$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped from
# 45182 (2023-09-13) to 45184 (2023-09-15) for every data row (2 through 407).
$ws.Range("C2:C407").Value = 45184
